# Applies the change described by the diff:
# - On the "Metadata" sheet, update the Date value and insert a new
#   "Jurisdiction" row (with an empty value) right after "Contact",
#   pushing Description/Purpose/Copyright/Immutable down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date property value (row 8, column B)
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# Insert a new row before current row 11 ("Description"), so everything
# from "Description" down shifts from row 11 -> row 12, etc.
$ws.Rows.Item(11).Insert()

# Match the formatting used by the other data rows (row 10, "Contact").
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row 11 with the Jurisdiction property
# (its value is left blank, same as the existing blank property rows).
$ws.Range("A11").Value = "Jurisdiction"
